$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Replace the body text of the "Installation" paragraph (the last
#    paragraph in the document) with the new, expanded wording.
# ------------------------------------------------------------------
$installPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $installPara.Range
$r.End = $r.End - 1
$r.Text = ""
$installPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $installPara.Range
$r.End = $r.End - 1
$r.InsertAfter("The software package can be installed in MATLAB by running setup.m, which will add the tools to the MATLAB path. The user can specify they`u{2019}re developing the toolbox further by adding the string `u{201C}develop`u{201D} as the second argument. This will create a directory Test that will be added within the main GeoData directory. This folder will not be added to the path.")

# ------------------------------------------------------------------
# Helper: append a new paragraph after the current last paragraph of
# the document, optionally bold, optionally with text.
# ------------------------------------------------------------------
function Add-Para($text, [bool]$bold) {
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $rr = $last.Range
    $rr.End = $rr.End - 1
    $rr.InsertParagraphAfter()
    $new = $d.Paragraphs($d.Paragraphs.Count)
    if ($text -ne $null -and $text -ne "") {
        $nr = $new.Range
        $nr.End = $nr.End - 1
        $nr.InsertAfter($text)
    }
    if ($bold) {
        $new.Range.Bold = 1
    }
    return $new
}

# empty paragraph
Add-Para "" $false | Out-Null

Add-Para "The python version can be installed using the setup tools package. To do this the user can type into the command line in the GeoData directory." $false | Out-Null

Add-Para "" $false | Out-Null

Add-Para "python setup.m" $false | Out-Null

Add-Para "" $false | Out-Null

Add-Para "To put get the development mode type " $false | Out-Null

Add-Para "" $false | Out-Null

Add-Para "python setup.m develop" $false | Out-Null

Add-Para "This will also place a directory called test in the main GeoData directory." $false | Out-Null

Add-Para "" $false | Out-Null

Add-Para "Format for h5 files" $true | Out-Null

Add-Para "" $true | Out-Null

Add-Para "The GeoData code base uses specifically formatted versions of h5 files. The purpose of these files is to easily read into both code bases. The format is set up so each group or data set in the base directory `u{2018}/`u{2019} will be a variable for the data set object that an instance of the class will represent.  " $false | Out-Null

Add-Para "" $false | Out-Null

$final = Add-Para "A data set will represent the content of the variable that is the name of the data set. The groups in the in the base directory of the h5 file will be read into Python as a dictionary where the names of the sub datasets will be keys and the data sets themselves will be the values. This is similar for MATLAB except that the groups will be read in as a struct with the dataset names representing the fieldnames and the datasets represent the values." $false

# ------------------------------------------------------------------
# 2. Move the hidden "_GoBack" bookmark from the "times ..." paragraph
#    to the end of the newly-added final paragraph (this also removes
#    it from its old location since a bookmark name is unique).
# ------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $final.Range) | Out-Null

Write-Output "done"
